# Updated cryptos list on Wed Mar  1 18:33:53 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "23.720.07"
$ws.Range("E2").Value = "  +1.07%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.656.29"
$ws.Range("E3").Value = "  +0.95%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - USDC
$ws.Range("E5").Value = "  +0.10%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'302.73"
$ws.Range("E6").Value = "  -0.19%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.3822"
$ws.Range("E7").Value = "  +0.54%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.06%  "

# Row 9 - OKB
$ws.Range("D9").Value = "'51.08"
$ws.Range("E9").Value = "  -1.98%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.08194"
$ws.Range("E10").Value = "  +0.27%  "

# Row 11 - Polygon
$ws.Range("D11").Value = "'1.231"
$ws.Range("E11").Value = "  -0.14%  "

# Row 12 - BinanceUSD
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  +0.02%  "

# Row 13 - Solana
$ws.Range("D13").Value = "'22.52"
$ws.Range("E13").Value = "  -0.02%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'6.451"
$ws.Range("E14").Value = "  -0.03%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "'7.437"
$ws.Range("E15").Value = "  +1.16%  "

# Row 16 - ShibaInu
$ws.Range("D16").Value = "'0.00001226"
$ws.Range("E16").Value = "  -0.87%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "1.655.18"
$ws.Range("E17").Value = "  +1.06%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "'97.65"
$ws.Range("E18").Value = "  +2.60%  "

# Row 19 - TRON
$ws.Range("D19").Value = "'0.07028"
$ws.Range("E19").Value = "  +1.03%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "'6.818"
$ws.Range("E20").Value = "  +3.70%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "'17.59"
$ws.Range("E21").Value = "  +0.33%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.02%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "'12.75"
$ws.Range("E23").Value = "  +1.99%  "

# Row 24 - WrappedBTC
$ws.Range("D24").Value = "23.731.31"
$ws.Range("E24").Value = "  +1.06%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "'2.500"
$ws.Range("E25").Value = "  -0.96%  "

# Row 26 - LidoDAOToken
$ws.Range("D26").Value = "'3.021"
$ws.Range("E26").Value = "  -1.34%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'21.24"
$ws.Range("E27").Value = "  +0.20%  "

# Row 28 - Monero
$ws.Range("D28").Value = "'154.54"
$ws.Range("E28").Value = "  +1.76%  "

# Row 29 - HuobiToken
$ws.Range("D29").Value = "'5.228"
$ws.Range("E29").Value = "  -0.81%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "'134.02"
$ws.Range("E30").Value = "  +0.54%  "

# Row 31 - WrappedliquidstakedEther2.0
$ws.Range("D31").Value = "1.838.67"
$ws.Range("E31").Value = "  +1.09%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'7.133"
$ws.Range("E32").Value = "  +8.45%  "

# Row 33 - WEMIXTOKEN
$ws.Range("D33").Value = "'2.243"
$ws.Range("E33").Value = "  +4.11%  "

# Row 34 - FraxShare
$ws.Range("D34").Value = "'11.96"
$ws.Range("E34").Value = "  +3.80%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  -3.22%  "

# Row 36 - VeChain
$ws.Range("D36").Value = "'0.02809"
$ws.Range("E36").Value = "  +1.37%  "

# Row 37 - Algorand
$ws.Range("D37").Value = "'0.2516"
$ws.Range("E37").Value = "  +0.19%  "

# Row 38 - was InternetComputer(DFINITY), now Stellar
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "'0.08803"
$ws.Range("E38").Value = "  +0.61%  "

# Row 39 - was Stellar, now InternetComputer(DFINITY)
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'6.091"
$ws.Range("E39").Value = "  +1.71%  "

# Row 40 - Hedera
$ws.Range("D40").Value = "'0.07003"
$ws.Range("E40").Value = "  -0.42%  "

# Row 41 - Aptos
$ws.Range("D41").Value = "'12.96"
$ws.Range("E41").Value = "  +5.64%  "

# Row 42 - TheSandbox
$ws.Range("D42").Value = "'0.6999"
$ws.Range("E42").Value = "  -0.66%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").Value = "'1.332"
$ws.Range("E43").Value = "  -1.37%  "

# Row 44 - EnergySwap
$ws.Range("D44").Value = "'16.03"
$ws.Range("E44").Value = "  +2.61%  "

# Row 45 - Decentraland
$ws.Range("D45").Value = "'0.6514"
$ws.Range("E45").Value = "  -0.39%  "

# Row 46 - Frax
$ws.Range("D46").Value = "'1.001"

# Row 47 - NEARProtocol
$ws.Range("D47").Value = "'2.302"
$ws.Range("E47").Value = "  +0.54%  "

# Row 48 - PancakeSwap
$ws.Range("D48").Value = "'3.964"
$ws.Range("E48").Value = "  +0.05%  "

# Row 49 - Cronos
$ws.Range("D49").Value = "'0.07902"
$ws.Range("E49").Value = "  -0.92%  "

# Row 50 - Quant
$ws.Range("D50").Value = "'128.20"
$ws.Range("E50").Value = "  -0.58%  "

# Row 51 - Flow
$ws.Range("D51").Value = "'1.180"
$ws.Range("E51").Value = "  -1.09%  "
